$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 19-24, same look & feel (styles) as the row above (row 18).
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C24").PasteSpecial(-4122)

# Row 19: Figure 5a
$ws.Range("A19").Value = "Figure 5a"
$ws.Range("B19").Value = "Bac2Feature_experiment/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb"
$ws.Range("C19").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb"

# Row 20: Figure 5b (shares File Path / GitHub URL with row 19)
$ws.Range("A20").Value = "Figure 5b"
$ws.Range("B20").Value = "Bac2Feature_experiment/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb"
$ws.Range("C20").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb"

# Row 21: Figure S2 vs homology-based prediction comparison
$ws.Range("A21").Value = "Figure S2"
$ws.Range("B21").Value = "Bac2Feature_experiment/scripts/09_cross_validation_suppl/092_homology/092_compare_homology_based_prediction.ipynb"
$ws.Range("C21").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/092_homology/092_compare_homology_based_prediction.ipynb"

# Row 22: Figure S3 vs taxonomic classifier comparison
$ws.Range("A22").Value = "Figure S3"
$ws.Range("B22").Value = "Bac2Feature_experiment/scripts/09_cross_validation_suppl/091_taxonomy/091_compare_taxonomic_classifier.ipynb"
$ws.Range("C22").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/091_taxonomy/091_compare_taxonomic_classifier.ipynb"

# Row 23: Figure S4 vs hsp method comparison
$ws.Range("A23").Value = "Figure S4"
$ws.Range("B23").Value = "Bac2Feature_experiment/scripts/09_cross_validation_suppl/093_phylogeny/093_compare_hsp_method.ipynb"
$ws.Range("C23").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/093_phylogeny/093_compare_hsp_method.ipynb"

# Row 24: Figure S7 clade-out cross validation
$ws.Range("A24").Value = "Figure S7"
$ws.Range("B24").Value = "Bac2Feature_experiment/scripts/09_cross_validation_suppl/094_clade_out/094_clade_out_cross_validation.ipynb"
$ws.Range("C24").Value = "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/094_clade_out/094_clade_out_cross_validation.ipynb"

# --- Hyperlinks for the new GitHub URL cells, added in the same order the
# author's file lists new relationship ids (rId16..rId21).
$ws.Hyperlinks.Add($ws.Range("C19"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb")
$ws.Hyperlinks.Add($ws.Range("C22"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/091_taxonomy/091_compare_taxonomic_classifier.ipynb")
$ws.Hyperlinks.Add($ws.Range("C23"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/093_phylogeny/093_compare_hsp_method.ipynb")
$ws.Hyperlinks.Add($ws.Range("C24"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/094_clade_out/094_clade_out_cross_validation.ipynb")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/08_application_lakesoilmicrobiome/08_application_lakesoilmicrobiome.ipynb")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://github.com/fuyo780/Bac2Feature_experiment/blob/main/scripts/09_cross_validation_suppl/092_homology/092_compare_homology_based_prediction.ipynb")

# Restore the body-row border/font styling the new hyperlinks wiped out.
$ws.Range("C18").Copy()
$ws.Range("C19:C24").PasteSpecial(-4122)

# --- Update rows 8 & 9: they used to point at the "Figure S2"/"Figure S3"
# codes (shared strings 7/8) with the old style (s=5); now they become new
# "Figure S5"/"Figure S6" codes with the plain body style (s=3), matching
# rows 6/7 above them. Done last so these two new shared strings land at
# the end of the table, matching the author's edit order.
$ws.Range("A6:C6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Figure S5"

$ws.Range("A6:C6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Figure S6"

# --- Cosmetic follow-up: widen column C to fit the new (longer) URLs, and
# leave the selection where the author's last edit left it.
$ws.Columns.Item(3).ColumnWidth = 144.15
$ws.Range("A20").Select()
